$wb = $excel.ActiveWorkbook

function Set-SheetRowsFlat {
    param($ws, $flat)
    for ($i = 0; $i -lt $flat.Count; $i += 6) {
        $r = $flat[$i]
        $ws.Cells.Item($r, 1).Value = $flat[$i + 1]
        $ws.Cells.Item($r, 2).Value = $flat[$i + 2]
        $ws.Cells.Item($r, 3).Value = $flat[$i + 3]
        $ws.Cells.Item($r, 4).Value = $flat[$i + 4]
        $ws.Cells.Item($r, 5).Value = $flat[$i + 5]
    }
}

$flat1 = @(
    112, "08:48:08", "09:35", "23_HERNANDEZ", 47, "LP1912", 113, "08:32:09", "09:35", "16_SANTA ANA", 63,
    "LP1912", 190, "10:59:49", "12:21", "26_HERNANDEZ", 82, "LP1912", 192, "12:21:08", "12:21", "16_SANTA ANA",
    0, "LP1912", 208, "12:59:47", "13:00", "16_SANTA ANA", 1, "LP1912", 209, "11:30:45", "13:01", "17_ROMERO",
    91, "LP1912", 210, "12:47:27", "13:02", "15_ABASTO", 15, "LP1912", 211, "12:21:08", "13:03", "14_ABASTO",
    42, "LP1912", 212, "12:47:27", "13:04", "23_HERNANDEZ", 17, "LP1912", 213, "12:59:47", "13:05", "23_HERNANDEZ",
    6, "LP1912", 214, "11:30:45", "13:06", "16_P MOR-SANTA ANA", 96, "LP1912", 215, "11:30:45", "13:07",
    "10_OLMOS", 97, "LP1912", 216, "12:21:08", "13:07", "16_P MOR-SANTA ANA", 46, "LP1912", 217, "12:21:08",
    "13:08", "10_OLMOS", 47, "LP1912", 218, "11:30:45", "13:13", "215D_EL PATO", 103, "LP1912", 219, "12:21:08",
    "13:14", "215D_EL PATO", 53, "LP1912", 220, "12:47:27", "13:14", "11_ETCHEVERRY", 27, "LP1912", 221,
    "11:56:55", "13:20", "26_HERNANDEZ", 84, "LP1912", 222, "11:30:45", "13:21", "26_HERNANDEZ", 111, "LP1912",
    223, "11:30:45", "13:25", "10_OLMOS", 115, "LP1912", 224, "11:30:45", "13:26", "15_ABASTO", 116, "LP1912",
    225, "11:30:45", "13:26", "14_ABASTO", 116, "LP1912", 226, "11:56:55", "13:27", "10_OLMOS", 91, "LP1912",
    227, "12:21:08", "13:27", "14_ABASTO", 66, "LP1912", 228, "12:21:08", "13:28", "10_OLMOS", 67, "LP1912",
    229, "12:47:27", "13:31", "10_OLMOS", 44, "LP1912", 230, "12:47:27", "13:32", "10_OLMOS", 45, "LP1912",
    231, "12:59:47", "13:33", "10_OLMOS", 34, "LP1912", 232, "11:56:55", "13:36", "15_ABASTO", 100, "LP1912",
    233, "11:56:55", "13:46", "17_ROMERO", 110, "LP1912", 234, "11:56:55", "13:50", "215A_EL PATO", 114,
    "LP1912", 235, "12:59:47", "13:50", "11_ETCHEVERRY", 51, "LP1912", 236, "12:21:08", "13:51", "215A_EL PATO",
    90, "LP1912", 237, "11:56:55", "13:55", "225_GOMEZ", 119, "LP1912", 238, "12:21:08", "13:56", "225_GOMEZ",
    95, "LP1912", 239, "12:59:47", "13:56", "16_P MOR-167 Y 521", 57, "LP1912", 240, "12:47:27", "13:58",
    "16_P MOR-167 Y 521", 71, "LP1912", 241, "12:21:08", "14:00", "16_P MOR-167 Y 521", 99, "LP1912", 242,
    "12:21:08", "14:04", "17_ROMERO", 103, "LP1912", 243, "12:21:08", "14:08", "23_HERNANDEZ", 107, "LP1912",
    244, "12:59:47", "14:11", "23_HERNANDEZ", 72, "LP1912", 245, "12:47:27", "14:16", "27_EL RETIRO", 89,
    "LP1912", 246, "12:21:08", "14:17", "27_EL RETIRO", 116, "LP1912", 247, "12:59:47", "14:19", "215C_EL PATO",
    80, "LP1912", 248, "12:21:08", "14:20", "215C_EL PATO", 119, "LP1912", 249, "12:47:27", "14:21", "26_HERNANDEZ",
    94, "LP1912", 250, "12:47:27", "14:45", "14_ABASTO", 118, "LP1912", 251, "12:59:47", "14:56", "16_P MOR-SANTA ANA",
    117, "LP1912", 252, "12:59:47", "14:58", "215B_EL PATO", 119, "LP1912"
)

$flat2 = @(
    33, "12:59:47", "14:19", "215C_EL PATO", 80, "LP1912", 34, "12:21:08", "14:20", "215C_EL PATO", 119,
    "LP1912", 35, "12:59:47", "14:58", "215B_EL PATO", 119, "LP1912"
)

$flat3 = @(
    44, "12:59:47", "14:52", "215D_LA PLATA", 113, "L6203"
)

$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 12:59:47"
$ws1.Cells.Item(3,1).Value = "Total filas: 247"
Set-SheetRowsFlat $ws1 $flat1

$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 12:59:47"
$ws2.Cells.Item(3,1).Value = "Total filas: 30"
Set-SheetRowsFlat $ws2 $flat2

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 12:59:47"
$ws3.Cells.Item(3,1).Value = "Total filas: 39"
Set-SheetRowsFlat $ws3 $flat3
